$wb = $excel.ActiveWorkbook

# 1) Rename "work_status" -> "activity" (content unchanged)
$activitySheet = $wb.Worksheets.Item("work_status")
$activitySheet.Name = "activity"

# 2) Insert three new lookup sheets after "activity": format, occupation, review_status
$formatSheet = $wb.Worksheets.Add($null, $activitySheet)
$formatSheet.Name = "format"

$occupationSheet = $wb.Worksheets.Add($null, $formatSheet)
$occupationSheet.Name = "occupation"

$reviewStatusSheet = $wb.Worksheets.Add($null, $occupationSheet)
$reviewStatusSheet.Name = "review_status"

# --- Pre-seed the 11 new shared strings in the exact order they were
# first introduced by the original import so the resulting sharedStrings
# table lines up with the source workbook (удаленка, офис, гибрид, полная,
# неполный день, стажировка, Не выбрано, На рассмотрении, Отправлено
# тестовое, Отказ, Назначено собеседование). We stash them in a scratch
# column on each new sheet, then clear the scratch cells once the real
# data (below) has reused those same shared-string slots.
$formatSheet.Range("D1").Value = "удаленка"
$formatSheet.Range("D2").Value = "офис"
$formatSheet.Range("D3").Value = "гибрид"

$occupationSheet.Range("D1").Value = "полная"
$occupationSheet.Range("D2").Value = "неполный день"
$occupationSheet.Range("D3").Value = "стажировка"

$reviewStatusSheet.Range("D1").Value = "Не выбрано"
$reviewStatusSheet.Range("D2").Value = "На рассмотрении"
$reviewStatusSheet.Range("D3").Value = "Отправлено тестовое"
$reviewStatusSheet.Range("D4").Value = "Отказ"
$reviewStatusSheet.Range("D5").Value = "Назначено собеседование"

# 3) "format" sheet data: id/name pairs
$formatSheet.Range("A1").Value = "id"
$formatSheet.Range("B1").Value = "name"
$formatSheet.Range("A2").Value = 1
$formatSheet.Range("B2").Value = "офис"
$formatSheet.Range("A3").Value = 2
$formatSheet.Range("B3").Value = "удаленка"
$formatSheet.Range("A4").Value = 3
$formatSheet.Range("B4").Value = "гибрид"

# 4) "occupation" sheet data: id/name pairs
$occupationSheet.Range("A1").Value = "id"
$occupationSheet.Range("B1").Value = "name"
$occupationSheet.Range("A2").Value = 1
$occupationSheet.Range("B2").Value = "полная"
$occupationSheet.Range("A3").Value = 2
$occupationSheet.Range("B3").Value = "неполный день"
$occupationSheet.Range("A4").Value = 3
$occupationSheet.Range("B4").Value = "стажировка"
$occupationSheet.Columns.Item(2).AutoFit()

# 5) "review_status" sheet data: id/name pairs
$reviewStatusSheet.Range("A1").Value = "id"
$reviewStatusSheet.Range("B1").Value = "name"
$reviewStatusSheet.Range("A2").Value = 1
$reviewStatusSheet.Range("B2").Value = "Не выбрано"
$reviewStatusSheet.Range("A3").Value = 2
$reviewStatusSheet.Range("B3").Value = "На рассмотрении"
$reviewStatusSheet.Range("A4").Value = 3
$reviewStatusSheet.Range("B4").Value = "Отправлено тестовое"
$reviewStatusSheet.Range("A5").Value = 4
$reviewStatusSheet.Range("B5").Value = "Назначено собеседование"
$reviewStatusSheet.Range("A6").Value = 5
$reviewStatusSheet.Range("B6").Value = "Отказ"

# Clear the scratch column used to order the shared-string table
$formatSheet.Range("D1:D3").Clear()
$occupationSheet.Range("D1:D3").Clear()
$reviewStatusSheet.Range("D1:D5").Clear()

# 6) Selections on each new sheet (matches source workbook)
$formatSheet.Range("C3").Select()
$occupationSheet.Range("A4").Select()
$reviewStatusSheet.Range("E5").Select()

# 7) The last-added sheet ("review_status") becomes the active tab
$reviewStatusSheet.Select()
$reviewStatusSheet.Range("E5").Select()
